# Generate Report for Archive
#
# 1. Update the localization status text from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2/F2, zh-cn!C2,
#    de-de!C2 all share the same string value).
# 2. Narrow the "Status" columns (Overview columns E/F, zh-cn/de-de column C)
#    to match the new, shorter text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ([string]$cell.Text -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
